$wb = $excel.ActiveWorkbook

# Sheet index 3: Restricciones_del_follower
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = "5.35 - 2x_1 + y_1 - y_2"
$ws3.Range("B2").Value = "'-2.8499999999999996"
$ws3.Range("D2").Value = "'0.73"
$ws3.Range("E2").Value = "'4.699999999999999"
$ws3.Range("F2").Value = "'3.9000000000000004"
$ws3.Range("A3").Value = "2.1499999999999932 + x_1 - 3x_2 + y_2"
$ws3.Range("B3").Value = "'-4.149999999999993"
$ws3.Range("D3").Value = "'0.59"
$ws3.Range("E3").Value = "'5.2"
$ws3.Range("F3").Value = "'2.7"
$ws3.Range("A4").Value = "104.95 - y_1"
$ws3.Range("B4").Value = "'-104.95"
$ws3.Range("D4").Value = "'0.87"
$ws3.Range("E4").Value = "'9.9"
$ws3.Range("F4").Value = "'2.8000000000000003"
$ws3.Range("A5").Value = "-3.5999999999999996 - y_2"
$ws3.Range("B5").Value = "'-3.5999999999999996"
$ws3.Range("D5").Value = "'0.08"
$ws3.Range("E5").Value = "'5.300000000000001"
$ws3.Range("F5").Value = "'7.9"

# Sheet index 4: Punto_modificado
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A2").Value = "'53.35"
$ws4.Range("B2").Value = "'19.7"
$ws4.Range("C2").Value = "'104.95"
$ws4.Range("D2").Value = "'3.5999999999999996"

# Sheet index 5: Vector_bf
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("A2").Value = "'4.14"
$ws5.Range("A3").Value = "'-0.78"

# Sheet index 6: Vector_BF
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("A2").Value = "'6.199999999999998"
$ws6.Range("A3").Value = "'14.600000000000001"
$ws6.Range("A4").Value = "'4.700000000000001"
$ws6.Range("A5").Value = "'4.8"
